$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first data row (old row 2, Datetime 45847.5625).
# This shifts every subsequent row up by one, so what used to be row 3
# becomes row 2, etc. Column A (Datetime) and column B (Real_Close)
# end up exactly where we want them after this shift.
$ws.Rows.Item(2).Delete()

# After the shift, column C (Trad_Prediction) must now hold the previous
# row's Real_Close value (i.e. yesterday's actual close used as the naive
# "traditional" prediction), and column D (AI_Prediction) must hold the
# newly generated AI prediction values. Build both arrays (rows 2-35) and
# write them back in one shot.

$tradPrediction = @(
  222.8899993896484, 222.3699951171875, 222.3999938964844, 222.6132965087891,
  222.4100036621094, 222.5807037353516, 222.4700012207031, 220.2599945068359,
  221.6347045898438, 222.0850067138672, 222.2899932861328, 222.6000061035156,
  222.0599975585938, 222.2799987792969, 224.3755950927734, 224.5299987792969,
  225.2601013183594, 225.5800018310547, 226.1000061035156, 225.7899932861328,
  224.9900054931641, 225.3056030273438, 225.8800048828125, 225.9949951171875,
  225.7550048828125, 225.4499969482422, 225.33349609375,   225.6300048828125,
  226.8600006103516, 226.6799926757812, 226.6000061035156, 225.9600067138672,
  226.8350067138672, 226.2969970703125
)

$aiPrediction = @(
  220.414228649649,  218.5475256756957, 228.1466958188372, 219.1970015261506,
  224.4246889943144, 226.5496740171144, 221.1789382914801, 220.6714637324821,
  225.8243371186124, 226.4673486776794, 226.5373763328199, 221.6177267275935,
  229.6165557793527, 220.8687229025682, 215.9071247420279, 220.7812943618412,
  218.7576481689501, 225.9436109821324, 230.3625239802085, 223.5342941934025,
  220.1100481040905, 225.8716682267265, 230.8914095450136, 222.8642131712212,
  224.6887210349146, 230.798391201034,  229.1735564823941, 232.3837217361472,
  224.9499663317515, 220.0809253772704, 222.732034681477,  222.067196652894,
  226.5572781674514, 231.4081065986274
)

$rowCount = $tradPrediction.Length
for ($i = 0; $i -lt $rowCount; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 3).Value = $tradPrediction[$i]
  $ws.Cells.Item($row, 4).Value = $aiPrediction[$i]
}
